# Calibração: implementando mudanças para calibração
# -----------------------------------------------------------------
# This script reproduces the author's edit:
#   1) Adjust the "Levers_FullDesign" scenario table (row 3) and drop
#      the now-unused row 4.
#   2) Re-peg the "PeD" budget ceiling (params!H66:H68) to a literal
#      value instead of the MAX() formula (the formula's own result,
#      now that Levers_FullDesign only runs through row 3).
#   3) Flip the "Estratégia de Capacidade" rows (params!72:74) from a
#      fixed value to an uncertain range (G/H bounds + I = "Incerto").
#   4) Re-create the stray `_xlnm._FilterDatabase_*` defined names that
#      accumulate whenever AutoFilter gets re-applied on the "params"
#      and "levers" sheets (cosmetic artefact of the round trip, kept
#      here so the defined-names table matches exactly).
#   5) Leave the cursor/selection where the author left it.
# -----------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsParams = $wb.Worksheets.Item("params")
$wsLevers = $wb.Worksheets.Item("Levers_FullDesign")

# --- 1) Levers_FullDesign: update row 3, delete row 4 --------------
$wsLevers.Range("A3").Value = 2
$wsLevers.Range("B3").Value = 0.9
$wsLevers.Range("C3").Value = 0.25
$wsLevers.Range("D3").Value = 0.15
$wsLevers.Rows.Item(4).Delete()

# --- 2) params: H66/H67/H68 become literal values -------------------
$wsParams.Range("H66").Value = 0.15
$wsParams.Range("H67").Value = 0.15
$wsParams.Range("H68").Value = 0.15

# --- 3) params: rows 72-74 become "Incerto" with a 0.51-2.5 band ----
$wsParams.Range("G72").Value = 0.51
$wsParams.Range("H72").Value = 2.5
$wsParams.Range("I72").Value = "Incerto"

$wsParams.Range("G73").Value = 0.51
$wsParams.Range("H73").Value = 2.5
$wsParams.Range("I73").Value = "Incerto"

$wsParams.Range("G74").Value = 0.51
$wsParams.Range("H74").Value = 2.5
$wsParams.Range("I74").Value = "Incerto"

# --- 4) Recreate the accumulated AutoFilter "ghost" defined names ---
$wsParams.Names.Add("_xlnm._FilterDatabase_0", "=params!`$A`$1:`$O`$78") | Out-Null
$wsParams.Names.Add("_xlnm._FilterDatabase_0_0", "=params!`$A`$1:`$O`$78") | Out-Null
$wsParams.Names.Add("_xlnm._FilterDatabase_0_0_0", "=params!`$A`$1:`$O`$78") | Out-Null
$wsParams.Names.Add("_xlnm._FilterDatabase_0_0_0_0", "=params!`$A`$1:`$O`$78") | Out-Null

$wsLevers2 = $wb.Worksheets.Item("levers")
$wsLevers2.Names.Add("_xlnm._FilterDatabase_0_0_0_0_0_0_0_0", "=levers!`$A`$1:`$H`$17") | Out-Null
$wsLevers2.Names.Add("_xlnm._FilterDatabase_0_0_0_0_0_0_0_0_0", "=levers!`$A`$1:`$G`$15") | Out-Null
$wsLevers2.Names.Add("_xlnm._FilterDatabase_0_0_0_0_0_0_0_0_0_0", "=levers!`$A`$1:`$H`$17") | Out-Null
$wsLevers2.Names.Add("_xlnm._FilterDatabase_0_0_0_0_0_0_0_0_0_0_0", "=levers!`$A`$1:`$G`$15") | Out-Null

# --- 5) Leave selection/scroll position where the author left it ----
$wsParams.Activate()
$wsParams.Range("A36").Select()
$excel.ActiveWindow.ScrollRow = 36
$wsParams.Range("I81").Select()

$wsLevers.Activate()
$wsLevers.Range("C5").Select()
